$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 1. First paragraph: insert a new bold run "产品名称" before the existing run
#    and change remaining text to "：神秘香料臻品印度奶茶"
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
# Collapse to the start of paragraph 1 to insert new bold text
$insertRange = $d.Range($r1.Start, $r1.Start)
$insertRange.Text = "产品名称"
$insertRange.Bold = $true

ReplaceText "产品名称：神秘香料高级柴茶" "：神秘香料臻品印度奶茶"

# 2. "主要特点：" -> "主要功能：" (and bold on)
$d.Content.Find.Execute("主要特点：", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
if ($d.Content.Find.Found) {
}
ReplaceText "主要特点：" "主要功能："

ReplaceText "正宗混合" "正宗配方"
ReplaceText "：我们的柴是优质黑茶叶的和谐混合，也是地香料的标志性选择，包括肉桂、豆瓜、丁香、姜和黑胡椒。" "：我们的奶茶选用优质黑茶，与肉桂、豆蔻、丁香、姜和黑胡椒等多种特色香料完美融合。"
ReplaceText "健康增强成分" "成分更加健康"
ReplaceText "：神秘香料柴茶中的每个成分都是出于自然健康益处而选择的。" "：神秘香料奶茶臻选自然原料，有利于健康。"
ReplaceText "浓郁的香气和味道" "香气浓郁、口味醇厚"
ReplaceText "：温暖，辣味和深，令人振奋的味道，我们的柴使它成为完美的饮料，开始你的一天或放松在晚上。" "：我们的奶茶气味温辛、口感醇厚，提神醒脑，是开启美好一天或晚上放松身心的完美饮品。"
ReplaceText "多才多艺的酿造选项" "多元化的烹制选项"
ReplaceText "：无论你喜欢你的柴热，作为一个令人耳目一新的冰茶，或作为奶油拿铁，我们的混合是多才多艺的，以满足任何偏好。" "：无论你是喜欢温热的奶茶，还是令人耳目一新的冰茶，或者是奶油拿铁，这款产品可以满足任何偏好。"
ReplaceText "可持续来源" "原料可持续"
ReplaceText "：致力于可持续性，我们从小型农场采购我们的成分，实践有机农业，不仅确保最好的品质，而且确保我们星球的福利。" "：我们注重可持续性，从小型农场采购原料，坚持有机农业，不仅能够确保极佳品质，而且可以确保对我们的星球有益。"
ReplaceText "优雅的包装" "包装精致"
ReplaceText "：神秘的香料柴茶是设计精美的生态友好包装，使其成为茶爱好者的理想礼物或豪华的礼物为自己。" "：神秘香料印度奶茶设计精美，采用生态友好的包装方式，因此是送给茶叶爱好者的理想礼物，也是送给自己的奢华之选。"
ReplaceText "：我们站在产品后面，提供满意保证。" "：我们为产品背书，提供令人满意的保证。"
ReplaceText "理想的选择" "适用人群"
ReplaceText "：茶爱好者、有健康意识的个人、温暖、辛辣的饮料爱好者，以及任何希望探索传统印度柴的丰富口味的人。" "：茶叶爱好者、注重健康的个人、喜欢温辛饮料的群体，以及希望品尝传统印度奶茶丰富口感的人群。"
